# Fix a typo ("West or your office" -> "West of your office") in the
# second paragraph of the "TextBox 7" shape on slide 11.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shp = $s.Shapes.Item(2)

$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(2, 1)

# "West or " starts right after the 15-character prefix "You live North ".
$sub = $para.Characters(16, 8)
$sub.Text = "West of "
